# Finalized definition of virtual environment
# Update calibration outputs across the AR / SETAR / GARCH / TARCH / AR_TARCH sheets.

$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws_AR = $wb.Worksheets.Item("AR")
$ws_AR.Range("B2").Value = -0.0000009876121723689516
$ws_AR.Range("B3").Value = 0.7779485607697052
$ws_AR.Range("B4").Value = 0.0001238913583820764
$ws_AR.Range("B5").Value = "[1.0, 0.1873333000678464, 0.15989578872708451, 0.11626070785454622, 0.10876975913107793, 0.350997585574195, 0.13713545332424257, 0.1272433161005007, 0.10332146253257464, 0.09095276999155319, 0.05123366092138888, 0.06277623973489015, 0.06821395031864988, 0.1061184329537056, 0.1343453511724172, 0.10333429991186058, 0.058335975817348576, 0.022203821518192685, 0.055624046857740396, 0.0794296460373645]"

# --- SETAR sheet ---
$ws_SETAR = $wb.Worksheets.Item("SETAR")
$ws_SETAR.Range("B3").Value = 0.9517684887459807
$ws_SETAR.Range("B4").Value = -0.004218802717399324
$ws_SETAR.Range("B5").Value = 0.6891555978228384
$ws_SETAR.Range("B6").Value = 0.00008714108966329977
$ws_SETAR.Range("B7").Value = 0.004848092790940782
$ws_SETAR.Range("B8").Value = 0.6235916877011048
$ws_SETAR.Range("B9").Value = 0.00006910064683416391
$ws_SETAR.Range("B10").Value = "[0.9999999999999999, 0.19274740975337826, 0.17863288049736553, 0.1725629772628786, 0.13167441178245487, 0.22482260932903417, 0.13373839791816197, 0.1206297423717175, 0.08697401093940876, 0.1131989678200537, 0.0898575455680258, 0.07125574558675808, 0.060366348120504294, 0.08193788220872573, 0.11410330540668653, 0.07908147662351217, 0.06846799796251873, 0.05606624707904407, 0.04041045068113377, 0.0538187576527082]"

# --- GARCH sheet ---
$ws_GARCH = $wb.Worksheets.Item("GARCH")
$ws_GARCH.Range("B2").Value = -0.0001024637953301668
$ws_GARCH.Range("B3").Value = 0.00000278780529691291
$ws_GARCH.Range("B4").Value = 0.09999999791540318
$ws_GARCH.Range("B5").Value = 0.8799999816555479
$ws_GARCH.Range("B6").Value = "[0.9999999999999999, -0.014956569221195672, -0.037938405929444266, -0.058800682010841196, -0.08861236910342102, 0.2670971822762601, 0.0024866761280456693, 0.03625380556219252, -0.0037487363936883553, -0.0356862734468483, -0.06444497919471968, -0.05538384207485115, -0.03433465833814008, 0.005075150999292938, 0.049022571242987215, 0.024379307548419927, -0.03117608752834977, -0.08675870836231679, -0.04912994615634956, -0.0020935011942433816]"

# --- TARCH sheet ---
$ws_TARCH = $wb.Worksheets.Item("TARCH")
$ws_TARCH.Range("B2").Value = -0.0001210012476741032
$ws_TARCH.Range("B3").Value = 0.000002787805297361218
$ws_TARCH.Range("B4").Value = 0.09999999793148424
$ws_TARCH.Range("B5").Value = 0.8749999819004871
$ws_TARCH.Range("B6").Value = "[1.0, -0.018424165952097574, -0.04087257884018654, -0.06074400102566728, -0.09095607330997228, 0.2646456825111537, 0.0006926385550862907, 0.03515439485138993, -0.004883447409422598, -0.03660265808000995, -0.06476153087539963, -0.055719039971281456, -0.03478348039590044, 0.0048324013159589, 0.04867596919507801, 0.02423130534519844, -0.031391929720862376, -0.08669303558435171, -0.049253769337589674, -0.0023088554969375266]"
$ws_TARCH.Range("B7").Value = 0.009999979107990777

# --- AR_TARCH sheet ---
$ws_AR_TARCH = $wb.Worksheets.Item("AR_TARCH")
$ws_AR_TARCH.Range("B2").Value = -0.9572257340566337
$ws_AR_TARCH.Range("B3").Value = 0.000002478257337105272
$ws_AR_TARCH.Range("B4").Value = 0.1000002937345969
$ws_AR_TARCH.Range("B5").Value = 0.8749997353894955
$ws_AR_TARCH.Range("B6").Value = "[1.0, 0.7354938038808728, 0.6257175474692295, 0.5273102186086286, 0.45300221106459904, 0.39703638613147807, 0.35684372394216635, 0.29706171005758164, 0.27975934339783937, 0.22613080073863737, 0.19710911042603455, 0.15702737711325362, 0.1600549349511534, 0.1141718676524015, 0.12941979485867208, 0.09579204829345618, 0.11466631345530649, 0.04657724252489246, 0.06616796429967463, 0.06466293059771713]"
$ws_AR_TARCH.Range("B7").Value = 0.009999409117294467
$ws_AR_TARCH.Range("B9").Value = 0.7789996684348611
